# Update performance dashboard 2025-12-22 23:33
# Refreshes "Annual Return (%)" (column G) and "Trading Days" (column M)
# for every (Pattern, Model) row across all worksheets (Summary + the three
# per-pattern sheets).

$wb = $excel.ActiveWorkbook

# Map of "Pattern|Model" -> (new Annual Return %, new Trading Days)
$updates = @{
    "Pattern1-Pure Data|deepseek-v3"          = @("+236.65%", 4)
    "Pattern1-Pure Data|gemini-3-pro"         = @("+218.11%", 4)
    "Pattern1-Pure Data|gpt-5"                = @("+99.85%", 4)
    "Pattern1-Pure Data|llama-3.1-405b"       = @("+29.33%", 4)
    "Pattern1-Pure Data|qwen-72b"             = @("+3.12%", 4)
    "Pattern2-Data+Technical|deepseek-v3"     = @("+17.67%", 4)
    "Pattern2-Data+Technical|gemini-3-pro"    = @("+29.45%", 4)
    "Pattern2-Data+Technical|gpt-5"           = @("+14.60%", 4)
    "Pattern2-Data+Technical|llama-3.1-405b"  = @("+5.46%", 4)
    "Pattern2-Data+Technical|qwen-72b"        = @("+2.03%", 4)
    "Pattern3-Data+News|deepseek-v3"          = @("+11.50%", 4)
    "Pattern3-Data+News|gemini-3-pro"         = @("+118.08%", 4)
    "Pattern3-Data+News|gpt-5"                = @("+165.01%", 4)
    "Pattern3-Data+News|llama-3.1-405b"       = @("+4.46%", 4)
    "Pattern3-Data+News|qwen-72b"             = @("+36.30%", 4)
}

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    # Format column G as Text for the data rows up-front so the percentage
    # strings we write (e.g. "+236.65%") are kept as literal text instead of
    # being auto-parsed into numbers, matching the source workbook where
    # these are stored as plain strings.
    $gRange = $ws.Range("G2:G" + $rowCount)
    $gRange.NumberFormat = "@"

    for ($r = 2; $r -le $rowCount; $r++) {
        $pattern = $ws.Cells.Item($r, 1).Value()
        $model = $ws.Cells.Item($r, 2).Value()
        if ($null -eq $pattern -or $null -eq $model) {
            continue
        }
        $key = "$pattern|$model"
        if ($updates.ContainsKey($key)) {
            $vals = $updates[$key]
            $ws.Cells.Item($r, 7).Value = $vals[0]
            $ws.Cells.Item($r, 13).Value = $vals[1]
        }
    }

    # Restore the default cell style now that the text has been written so
    # no extra/visible formatting lingers on the updated cells.
    $gRange.Style = "Normal"
}
